$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.855.08'
$ws.Range("E2").Value = '  -0.72%  '

$ws.Range("D3").Value = '2.329.03'
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("E4").Value = '  +0.04%  '

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '302.80'
$c.ClearFormats()
$ws.Range("E5").Value = '  -0.30%  '

$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '93.89'
$c.ClearFormats()
$ws.Range("E6").Value = '  -3.87%  '

$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.501'
$c.ClearFormats()
$ws.Range("E7").Value = '  -1.03%  '

$ws.Range("E8").Value = '  +0.05%  '

$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.492'
$c.ClearFormats()
$ws.Range("E9").Value = '  -1.77%  '

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '33.95'
$c.ClearFormats()
$ws.Range("E10").Value = '  -4.55%  '

$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.0780'
$c.ClearFormats()
$ws.Range("E11").Value = '  -2.18%  '

$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '18.67'
$c.ClearFormats()
$ws.Range("E12").Value = '  -4.46%  '

$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '0.120'
$c.ClearFormats()
$ws.Range("E13").Value = '  +1.05%  '

$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '6.70'
$c.ClearFormats()
$ws.Range("E14").Value = '  -3.60%  '

$ws.Range("D15").Value = '2.690.87'
$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").Value = '2.350.55'
$ws.Range("E16").Value = '  +1.01%  '

$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '0.788'
$c.ClearFormats()
$ws.Range("E17").Value = '  +0.16%  '

$ws.Range("D18").Value = '42.788.81'
$ws.Range("E18").Value = '  -0.44%  '

$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '12.01'
$c.ClearFormats()
$ws.Range("E19").Value = '  -4.77%  '

$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '6.17'
$c.ClearFormats()
$ws.Range("E20").Value = '  +1.84%  '

$ws.Range("D21").Value = '0.0₃0887'
$ws.Range("E21").Value = '  -1.47%  '

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '67.83'
$c.ClearFormats()
$ws.Range("E22").Value = '  -0.04%  '

$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '234.99'
$c.ClearFormats()
$ws.Range("E23").Value = '  -0.69%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("E26").Value = '  -1.69%  '

$ws.Range("E27").Value = '  -1.86%  '

$ws.Range("E28").Value = '  +13.77%  '

$ws.Range("E29").Value = '  -0.54%  '

$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '31.10'
$c.ClearFormats()
$ws.Range("E30").Value = '  -6.40%  '

$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '0.0758'
$c.ClearFormats()
$ws.Range("E32").Value = '  +8.77%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '4.97'
$c.ClearFormats()
$ws.Range("E33").Value = '  -0.56%  '

$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '17.30'
$c.ClearFormats()
$ws.Range("E34").Value = '  -4.02%  '

$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '132.11'
$c.ClearFormats()
$ws.Range("E35").Value = '  -20.03%  '

$ws.Range("E36").Value = '  -1.14%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '1.81'
$c.ClearFormats()
$ws.Range("E37").Value = '  +2.60%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '4.33'
$c.ClearFormats()
$ws.Range("E38").Value = '  -4.46%  '

$ws.Range("E39").Value = '  -0.74%  '

$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '22.21'
$c.ClearFormats()
$ws.Range("E40").Value = '  +22.79%  '

$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '2.74'
$c.ClearFormats()
$ws.Range("E41").Value = '  -1.92%  '

$ws.Range("E42").Value = '  -1.43%  '

$ws.Range("D43").Value = '1.923.79'
$ws.Range("E43").Value = '  -3.39%  '

$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '0.0280'
$c.ClearFormats()
$ws.Range("E44").Value = '  +0.19%  '

$ws.Range("E45").Value = '  -5.34%  '

$ws.Range("E46").Value = '  -0.01%  '

$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '2.72'
$c.ClearFormats()
$ws.Range("E47").Value = '  -2.02%  '

$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '2.88'
$c.ClearFormats()
$ws.Range("E48").Value = '  -0.56%  '

$ws.Range("D49").Value = '2.558.63'
$ws.Range("E49").Value = '  +0.23%  '

$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '52.61'
$c.ClearFormats()
$ws.Range("E50").Value = '  -2.12%  '

$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '71.87'
$c.ClearFormats()
$ws.Range("E51").Value = '  -0.13%  '
